$wb = $excel.ActiveWorkbook

# --- Assets sheet: remove the three Captcha-related asset rows ---
# (Captcha_SiteKey/RPA021_MOLPAY_Captcha_SiteKey, Captcha_RuleId/RPA_Moon_Captcha_RuleId,
#  GCaptcha_RuleId/RPA_Moon_GCaptcha_RuleId) which occupy rows 6-8.
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Range("A6:A7").EntireRow.Select()
$wsAssets.Range("A6:A8").EntireRow.Delete()

# --- Settings sheet: add new e2e Air Asia OTP credential row ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
$wsSettings.Range("B8").Value = "RPA106_Air Asia_Email OTP"
$wsSettings.Range("A8").Value = "Cred_OTP"
$wsSettings.Range("C8").Value = "Cred email to get OTP"
$wsSettings.Range("C8").WrapText = $true
$wsSettings.Range("A8").Select()
